# Shaalmi_Tests.xlsx edit script
# Summary of changes:
#  - Browser sheet: switch browser keyword from FIREFOX to CHROME
#  - New "GO2URL" / "url" navigation step inserted as row 2 in the
#    Admin_SignIn_Tests and Login_Test sheets (parallel run support)
#  - Admin_TransportTab_Tests becomes the active/selected sheet

$wb = $excel.ActiveWorkbook

# --- 1. Browser sheet: FIREFOX -> CHROME -------------------------------
$wsBrowser = $wb.Worksheets.Item("Browser")
$wsBrowser.Range("A2").Value = "CHROME"

# --- 2. Admin_SignIn_Tests: insert GO2URL/url row at row 2 -------------
$wsSignIn = $wb.Worksheets.Item("Admin_SignIn_Tests")

# drop existing hyperlinks so stale refs don't linger after the shift
$wsSignIn.Range("A1").Hyperlinks.Delete()

$wsSignIn.Rows.Item(2).Insert()
$newRow = $wsSignIn.Range("A2:D2")
$newRow.ClearFormats()
$newRow.Borders.LineStyle = 1
$wsSignIn.Range("A2").Value = "GO2URL"
$wsSignIn.Range("D2").Value = "url"

$wsSignIn.Activate()
$wsSignIn.Range("A2:D2").Select()

# --- 3. Login_Test: insert GO2URL/url row at row 2 ---------------------
$wsLogin = $wb.Worksheets.Item("Login_Test")

# capture + drop existing hyperlinks (mailto links on D5 / D12 after shift)
$wsLogin.Range("A1").Hyperlinks.Delete()

$wsLogin.Rows.Item(2).Insert()
$newRow2 = $wsLogin.Range("A2:D2")
$newRow2.ClearFormats()
$newRow2.Borders.LineStyle = 1
$wsLogin.Range("A2").Value = "GO2URL"
$wsLogin.Range("D2").Value = "url"

# re-add the two mailto hyperlinks at their new (shifted down by one) cells
$wsLogin.Hyperlinks.Add($wsLogin.Range("D5"), "mailto:awais@gmail.com")
$wsLogin.Hyperlinks.Add($wsLogin.Range("D12"), "mailto:nauman@hotmail.com")

$wsLogin.Activate()
$wsLogin.Range("A2:D2").Select()

# --- 4. Admin_TransportTab_Tests becomes the active / selected sheet ---
$wsTransport = $wb.Worksheets.Item("Admin_TransportTab_Tests")
$wsTransport.Activate()
$wsTransport.Range("C4").Select()
